$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7937279.5
$ws.Range("I43").Value = 833.3333
$ws.Range("J43").Value = 13889614
$ws.Range("K43").Value = 833.3333
$ws.Range("L43").Value = 13889614
$ws.Range("M43").Value = -764.3333
$ws.Range("N43").Value = -13889752
$ws.Range("H53").Value = 1537.6666
$ws.Range("I53").Value = 2232.3333
$ws.Range("K53").Value = 2232.3333
$ws.Range("M53").Value = -1595.3333
$ws.Range("H62").Value = 3629.6667
$ws.Range("I62").Value = 4444.75
$ws.Range("K62").Value = 4444.75
$ws.Range("M62").Value = -3820.75
$ws.Range("H65").Value = 3629.6667
$ws.Range("I65").Value = 4444.75
$ws.Range("K65").Value = 22223.75
$ws.Range("M65").Value = -19103.75
$ws.Range("H132").Value = 5118.7666
$ws.Range("I132").Value = 2802.8635
$ws.Range("J132").Value = 11487.5
$ws.Range("K132").Value = 8408.5905
$ws.Range("L132").Value = 34462.5
$ws.Range("M132").Value = -5878.5905
$ws.Range("N132").Value = -39522.5
$ws.Range("H139").Value = 48878.184
$ws.Range("J139").Value = 51266
$ws.Range("L139").Value = 51266
$ws.Range("N139").Value = -61546
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3287.913
$ws.Range("I32").Value = 3129.318
$ws.Range("K32").Value = 3129.318
$ws.Range("M32").Value = -2842.318
$ws.Range("H61").Value = 1175.4
$ws.Range("I61").Value = 881.5714
$ws.Range("J61").Value = 1861
$ws.Range("K61").Value = 881.5714
$ws.Range("L61").Value = 1861
$ws.Range("M61").Value = -669.5714
$ws.Range("N61").Value = -2285
$ws.Range("H122").Value = 1581.1666
$ws.Range("I122").Value = 1496.75
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 4490.25
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -2040.25
$ws.Range("N122").Value = -10150
$ws.Range("H136").Value = 1175.4
$ws.Range("I136").Value = 881.5714
$ws.Range("J136").Value = 1861
$ws.Range("K136").Value = 2644.7142
$ws.Range("L136").Value = 5583
$ws.Range("M136").Value = -94.71420000000035
$ws.Range("N136").Value = -10683
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3286.25
$ws.Range("I86").Value = 3476.238
$ws.Range("K86").Value = 3476.238
$ws.Range("M86").Value = -2353.238
$ws.Range("H89").Value = 3286.25
$ws.Range("I89").Value = 3476.238
$ws.Range("K89").Value = 17381.19
$ws.Range("M89").Value = -11765.19
$ws.Range("H105").Value = 111113230
$ws.Range("I105").Value = 125002170
$ws.Range("K105").Value = 125002170
$ws.Range("M105").Value = -125000423
$ws.Range("H134").Value = 7272.4707
$ws.Range("I134").Value = 1375.5333
$ws.Range("K134").Value = 4126.5999
$ws.Range("M134").Value = -1591.5999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71429730
$ws.Range("I16").Value = 111112150
$ws.Range("J16").Value = 1378.6
$ws.Range("K16").Value = 111112150
$ws.Range("L16").Value = 1378.6
$ws.Range("M16").Value = -111111863
$ws.Range("N16").Value = -1952.6
$ws.Range("H22").Value = 345.42856
$ws.Range("I22").Value = 234.5
$ws.Range("J22").Value = 446.27274
$ws.Range("K22").Value = 234.5
$ws.Range("L22").Value = 446.27274
$ws.Range("M22").Value = 115.5
$ws.Range("N22").Value = -1146.27274
$ws.Range("H31").Value = 1273.8572
$ws.Range("I31").Value = 949.3
$ws.Range("J31").Value = 1568.909
$ws.Range("K31").Value = 949.3
$ws.Range("L31").Value = 1568.909
$ws.Range("M31").Value = -654.3
$ws.Range("N31").Value = -2158.909
$ws.Range("H34").Value = 1273.8572
$ws.Range("I34").Value = 949.3
$ws.Range("J34").Value = 1568.909
$ws.Range("K34").Value = 949.3
$ws.Range("L34").Value = 1568.909
$ws.Range("M34").Value = -747.3
$ws.Range("N34").Value = -1972.909
$ws.Range("H35").Value = 970
$ws.Range("I35").Value = 970
$ws.Range("K35").Value = 970
$ws.Range("M35").Value = -676
$ws.Range("H62").Value = 4547961.5
$ws.Range("I62").Value = 2565.2327
$ws.Range("K62").Value = 2565.2327
$ws.Range("M62").Value = -1941.2327
$ws.Range("H65").Value = 4547961.5
$ws.Range("I65").Value = 2565.2327
$ws.Range("K65").Value = 12826.1635
$ws.Range("M65").Value = -9706.163500000001
$ws.Range("H92").Value = 70320.2
$ws.Range("J92").Value = 70320.2
$ws.Range("L92").Value = 70320.2
$ws.Range("N92").Value = -75312.2
$ws.Range("H113").Value = 71429730
$ws.Range("I113").Value = 111112150
$ws.Range("J113").Value = 1378.6
$ws.Range("K113").Value = 111112150
$ws.Range("L113").Value = 1378.6
$ws.Range("M113").Value = -111109980
$ws.Range("N113").Value = -5718.6
$ws.Range("H122").Value = 776.625
$ws.Range("I122").Value = 719.8
$ws.Range("J122").Value = 871.3333
$ws.Range("K122").Value = 2159.4
$ws.Range("L122").Value = 2613.9999
$ws.Range("M122").Value = 290.6000000000004
$ws.Range("N122").Value = -7513.9999
$ws.Range("H134").Value = 1774.6666
$ws.Range("I134").Value = 1632.8889
$ws.Range("K134").Value = 4898.6667
$ws.Range("M134").Value = -2363.6667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 203.48
$ws.Range("I12").Value = 75.333336
$ws.Range("J12").Value = 321.76923
$ws.Range("K12").Value = 226.000008
$ws.Range("L12").Value = 965.30769
$ws.Range("M12").Value = -53.00000800000001
$ws.Range("N12").Value = -1311.30769
$ws.Range("H33").Value = 189
$ws.Range("I33").Value = 117.166664
$ws.Range("K33").Value = 702.999984
$ws.Range("M33").Value = -419.999984
$ws.Range("H98").Value = 895.6667
$ws.Range("I98").Value = 336.4
$ws.Range("K98").Value = 1009.2
$ws.Range("M98").Value = 488.8000000000001
$ws.Range("H122").Value = 841.86664
$ws.Range("I122").Value = 611.2857
$ws.Range("J122").Value = 1043.625
$ws.Range("K122").Value = 5501.571300000001
$ws.Range("L122").Value = 9392.625
$ws.Range("M122").Value = -3051.571300000001
$ws.Range("N122").Value = -14292.625
$ws.Range("H126").Value = 5159.8
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 5466.4443
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 16399.3329
$ws.Range("M126").Value = -2260
$ws.Range("N126").Value = -26279.3329
$ws.Range("H129").Value = 13890006
$ws.Range("J129").Value = 4387289
$ws.Range("L129").Value = 13161867
$ws.Range("N129").Value = -13171867
$ws.Range("H131").Value = 18183138
$ws.Range("J131").Value = 1421.8914
$ws.Range("L131").Value = 4265.674199999999
$ws.Range("N131").Value = -14345.6742
$ws.Range("H134").Value = 2777.5264
$ws.Range("I134").Value = 1522.8334
$ws.Range("J134").Value = 4928.4287
$ws.Range("K134").Value = 4568.5002
$ws.Range("L134").Value = 14785.2861
$ws.Range("M134").Value = 501.4997999999996
$ws.Range("N134").Value = -24925.2861
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21891.5
$ws.Range("J92").Value = 21891.5
$ws.Range("L92").Value = 21891.5
$ws.Range("N92").Value = -25635.5
$ws.Range("H97").Value = 1181.8
$ws.Range("I97").Value = 1181.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1181.8
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -685.8
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 1541.9286
$ws.Range("I102").Value = 1461.3
$ws.Range("J102").Value = 1743.5
$ws.Range("K102").Value = 1461.3
$ws.Range("L102").Value = 1743.5
$ws.Range("M102").Value = 160.7
$ws.Range("N102").Value = -4987.5
$ws.Range("H110").Value = 17621.2
$ws.Range("J110").Value = 17621.2
$ws.Range("L110").Value = 17621.2
$ws.Range("N110").Value = -25801.2
$ws.Range("H113").Value = 1879
$ws.Range("I113").Value = 1056.2222
$ws.Range("K113").Value = 1056.2222
$ws.Range("M113").Value = 1113.7778
$ws.Range("H132").Value = 2837.1365
$ws.Range("I132").Value = 2519.3572
$ws.Range("J132").Value = 3393.25
$ws.Range("K132").Value = 7558.071599999999
$ws.Range("L132").Value = 10179.75
$ws.Range("M132").Value = -5028.071599999999
$ws.Range("N132").Value = -15239.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1646.2727
$ws.Range("I7").Value = 1460.4
$ws.Range("K7").Value = 1460.4
$ws.Range("M7").Value = -1348.4
$ws.Range("H22").Value = 986.2727
$ws.Range("I22").Value = 516.6667
$ws.Range("J22").Value = 1162.375
$ws.Range("K22").Value = 516.6667
$ws.Range("L22").Value = 1162.375
$ws.Range("M22").Value = -221.6667
$ws.Range("N22").Value = -1752.375
$ws.Range("H27").Value = 986.2727
$ws.Range("I27").Value = 516.6667
$ws.Range("J27").Value = 1162.375
$ws.Range("K27").Value = 516.6667
$ws.Range("L27").Value = 1162.375
$ws.Range("M27").Value = -409.6667
$ws.Range("N27").Value = -1376.375
$ws.Range("H61").Value = 2974.2
$ws.Range("I61").Value = 2200.3333
$ws.Range("J61").Value = 4135
$ws.Range("K61").Value = 2200.3333
$ws.Range("L61").Value = 4135
$ws.Range("M61").Value = -1998.3333
$ws.Range("N61").Value = -4539
$ws.Range("H113").Value = 2974.2
$ws.Range("I113").Value = 2200.3333
$ws.Range("J113").Value = 4135
$ws.Range("K113").Value = 2200.3333
$ws.Range("L113").Value = 4135
$ws.Range("M113").Value = -30.33329999999978
$ws.Range("N113").Value = -8475
$ws.Range("H126").Value = 1646.2727
$ws.Range("I126").Value = 1460.4
$ws.Range("K126").Value = 4381.200000000001
$ws.Range("M126").Value = -1911.200000000001
$ws.Range("H132").Value = 79823.234
$ws.Range("I132").Value = 2866.5
$ws.Range("J132").Value = 145786.14
$ws.Range("K132").Value = 8599.5
$ws.Range("L132").Value = 437358.42
$ws.Range("M132").Value = -6069.5
$ws.Range("N132").Value = -442418.42
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 647.25
$ws.Range("I100").Value = 605.8
$ws.Range("J100").Value = 716.3333
$ws.Range("K100").Value = 1211.6
$ws.Range("L100").Value = 1432.6666
$ws.Range("M100").Value = -670.5999999999999
$ws.Range("N100").Value = -2514.6666
